# Update cryptocurrency price/volume data to match the latest GitHub Actions run.
# All touched cells hold plain text (prices use "." as thousands separators in
# some rows, e.g. "26.889.04", and percentage cells keep their padding
# whitespace), so we temporarily force Text formatting on the data range
# before writing the new values. Otherwise values that look like plain
# numbers (e.g. "309.42") would be auto-converted to numeric cells, which
# would not match the original inline-string cell representation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("B2:E51")
$origStyle = $dataRange.Style
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '26.889.04'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').Value = '1.843.51'
$ws.Range('E3').Value = '  +1.18%  '
$ws.Range('E4').Value = '  -0.10%  '
$ws.Range('D5').Value = '309.42'
$ws.Range('E5').Value = '  +1.38%  '
$ws.Range('E6').Value = '  -0.17%  '
$ws.Range('D7').Value = '0.4757'
$ws.Range('E7').Value = '  +2.40%  '
$ws.Range('D8').Value = '0.3666'
$ws.Range('E8').Value = '  +2.02%  '
$ws.Range('D9').Value = '0.07196'
$ws.Range('D10').Value = '0.9261'
$ws.Range('E10').Value = '  +3.14%  '
$ws.Range('D11').Value = '19.66'
$ws.Range('E11').Value = '  +1.81%  '
$ws.Range('D12').Value = '0.07694'
$ws.Range('E12').Value = '  -0.97%  '
$ws.Range('D13').Value = '1.898.39'
$ws.Range('E13').Value = '  +4.07%  '
$ws.Range('D14').Value = '5.311'
$ws.Range('E14').Value = '  +1.24%  '
$ws.Range('D15').Value = '6.397'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').Value = '88.71'
$ws.Range('E16').Value = '  +1.77%  '
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').Value = '0.000008639'
$ws.Range('E18').Value = '  +1.24%  '
$ws.Range('E19').Value = '  -0.11%  '
$ws.Range('B20').Value = 'Avalanche'
$ws.Range('C20').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D20').Value = '14.55'
$ws.Range('E20').Value = '  +3.13%  '
$ws.Range('B21').Value = 'WrappedBTC'
$ws.Range('C21').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D21').Value = '26.916.62'
$ws.Range('E21').Value = '  +1.12%  '
$ws.Range('D22').Value = '5.052'
$ws.Range('E22').Value = '  +0.83%  '
$ws.Range('D23').Value = '10.63'
$ws.Range('D24').Value = '1.918'
$ws.Range('E24').Value = '  +0.14%  '
$ws.Range('D25').Value = '152.30'
$ws.Range('E25').Value = '  +0.23%  '
$ws.Range('D26').Value = '18.14'
$ws.Range('E26').Value = '  +1.48%  '
$ws.Range('D27').Value = '2.000'
$ws.Range('E27').Value = '  +1.68%  '
$ws.Range('D28').Value = '114.13'
$ws.Range('E28').Value = '  +0.46%  '
$ws.Range('D29').Value = '4.921'
$ws.Range('E29').Value = '  +2.60%  '
$ws.Range('D30').Value = '0.08867'
$ws.Range('E30').Value = '  +0.86%  '
$ws.Range('D31').Value = '3.298'
$ws.Range('E31').Value = '  +4.99%  '
$ws.Range('B32').Value = 'ARBITRUM'
$ws.Range('C32').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D32').Value = '1.172'
$ws.Range('E32').Value = '  +4.01%  '
$ws.Range('B33').Value = 'ImmutableX'
$ws.Range('C33').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D33').Value = '0.7470'
$ws.Range('E33').Value = '  +2.49%  '
$ws.Range('D34').Value = '4.477'
$ws.Range('E34').Value = '  +1.09%  '
$ws.Range('D35').Value = '2.728'
$ws.Range('E35').Value = '  +0.23%  '
$ws.Range('D36').Value = '1.092'
$ws.Range('E36').Value = '  +1.67%  '
$ws.Range('D37').Value = '0.01953'
$ws.Range('E37').Value = '  +1.73%  '
$ws.Range('D38').Value = '0.05257'
$ws.Range('E38').Value = '  +3.21%  '
$ws.Range('D39').Value = '2.968'
$ws.Range('D40').Value = '0.5200'
$ws.Range('E40').Value = '  +3.43%  '
$ws.Range('D41').Value = '6.948'
$ws.Range('E41').Value = '  +1.41%  '
$ws.Range('D42').Value = '0.1509'
$ws.Range('E42').Value = '  +1.30%  '
$ws.Range('D43').Value = '8.202'
$ws.Range('E43').Value = '  +3.30%  '
$ws.Range('D44').Value = '10.55'
$ws.Range('E44').Value = '  +5.84%  '
$ws.Range('D45').Value = '0.4724'
$ws.Range('E45').Value = '  +1.90%  '
$ws.Range('D46').Value = '1.007'
$ws.Range('E46').Value = '  -0.17%  '
$ws.Range('D47').Value = '101.32'
$ws.Range('E47').Value = '  +3.23%  '
$ws.Range('D48').Value = '1.599'
$ws.Range('E48').Value = '  +3.02%  '
$ws.Range('D49').Value = '65.75'
$ws.Range('E49').Value = '  +3.36%  '
$ws.Range('D50').Value = '0.06019'
$ws.Range('E50').Value = '  +0.68%  '
$ws.Range('D51').Value = '0.8847'
$ws.Range('E51').Value = '  +4.18%  '

$dataRange.Style = $origStyle
